$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.335.91"
$ws.Range("E2").Value = "  +0.08%  "

$ws.Range("D3").Value = "3.498.64"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'589.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").Value = "'134.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "'7.70"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.05%  "

$ws.Range("D10").Value = "'0.125"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.15%  "

$ws.Range("E11").Value = "  +2.72%  "

$ws.Range("D12").Value = "4.093.23"
$ws.Range("E12").Value = "  +0.24%  "

$ws.Range("E13").Value = "  +0.68%  "

$ws.Range("E14").Value = "  -0.20%  "

$ws.Range("D15").Value = "3.498.35"
$ws.Range("E15").Value = "  +0.20%  "

$ws.Range("D16").Value = "64.274.84"
$ws.Range("E16").Value = "  -0.04%  "

$ws.Range("D17").Value = "'25.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.90%  "

$ws.Range("D18").Value = "'10.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.11%  "

$ws.Range("E19").Value = "  +0.53%  "

$ws.Range("D20").Value = "'13.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.27%  "

$ws.Range("D21").Value = "'386.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.33%  "

$ws.Range("D22").Value = "'0.580"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").Value = "3.637.20"
$ws.Range("E23").Value = "  +0.16%  "

$ws.Range("D24").Value = "'74.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.20%  "

$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("E26").Value = "  -0.37%  "

$ws.Range("D27").Value = "'0.0000116"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.90%  "

$ws.Range("D28").Value = "'0.992"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.74%  "

$ws.Range("D29").Value = "'7.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.53%  "

$ws.Range("D31").Value = "'1.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.31%  "

$ws.Range("D32").Value = "'8.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.23%  "

$ws.Range("E33").Value = "  +4.13%  "

$ws.Range("D34").Value = "3.526.66"
$ws.Range("E34").Value = "  +0.45%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("D36").Value = "'23.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.57%  "

$ws.Range("D37").Value = "'5.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.97%  "

$ws.Range("E38").Value = "  +0.65%  "

$ws.Range("E39").Value = "  +0.40%  "

$ws.Range("D40").Value = "'165.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.03%  "

$ws.Range("D41").Value = "'0.0785"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.70%  "

$ws.Range("D42").Value = "'0.808"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.44%  "

$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("E44").Value = "  +0.32%  "

$ws.Range("D45").Value = "'24.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.96%  "

$ws.Range("D46").Value = "'1.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.95%  "

$ws.Range("D47").Value = "'1.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.75%  "

$ws.Range("D48").Value = "2.444.39"
$ws.Range("E48").Value = "  -1.16%  "

$ws.Range("D49").Value = "'6.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.69%  "

$ws.Range("D50").Value = "'0.915"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.12%  "

$ws.Range("E51").Value = "  -0.42%  "
